# Auto-generated edit script: updates cached numeric values in H:N columns
# across all 8 sheets per the scheduled-runner data refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 182270.64
$ws.Range("I4").Value = 182270.64
$ws.Range("K4").Value = 182270.64
$ws.Range("M4").Value = -182156.64
$ws.Range("H55").Value = 272.18182
$ws.Range("I55").Value = 133.57143
$ws.Range("J55").Value = 514.75
$ws.Range("K55").Value = 133.57143
$ws.Range("L55").Value = 514.75
$ws.Range("M55").Value = 80.42857000000001
$ws.Range("N55").Value = -942.75
$ws.Range("H98").Value = 4016.6538
$ws.Range("I98").Value = 3629.3333
$ws.Range("K98").Value = 3629.3333
$ws.Range("M98").Value = -2131.3333
$ws.Range("H116").Value = 4500
$ws.Range("I116").Value = 4500
$ws.Range("K116").Value = 4500
$ws.Range("M116").Value = -1058
$ws.Range("H122").Value = 4016.6538
$ws.Range("I122").Value = 3629.3333
$ws.Range("K122").Value = 10887.9999
$ws.Range("M122").Value = -8437.999899999999
$ws.Range("H131").Value = 596006.6
$ws.Range("I131").Value = 722765.4
$ws.Range("J131").Value = 4466
$ws.Range("K131").Value = 2168296.2
$ws.Range("L131").Value = 13398
$ws.Range("M131").Value = -2163256.2
$ws.Range("N131").Value = -23478
$ws.Range("H132").Value = 6406.08
$ws.Range("I132").Value = 6406.08
$ws.Range("K132").Value = 19218.24
$ws.Range("M132").Value = -16688.24
$ws.Range("H135").Value = 604.3333
$ws.Range("J135").Value = 2499
$ws.Range("L135").Value = 22491
$ws.Range("N135").Value = -27561
$ws.Range("H137").Value = 2841.2693
$ws.Range("I137").Value = 1701
$ws.Range("K137").Value = 5103
$ws.Range("M137").Value = -2553
$ws.Range("H138").Value = 2691.6365
$ws.Range("I138").Value = 2037.25
$ws.Range("J138").Value = 2757.075
$ws.Range("K138").Value = 6111.75
$ws.Range("L138").Value = 8271.224999999999
$ws.Range("M138").Value = -971.75
$ws.Range("N138").Value = -18551.225
$ws.Range("H141").Value = 2849.7778
$ws.Range("I141").Value = 2755.6875
$ws.Range("K141").Value = 8267.0625
$ws.Range("M141").Value = -3087.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1661613
$ws.Range("I32").Value = 824860.4
$ws.Range("K32").Value = 824860.4
$ws.Range("M32").Value = -824573.4
$ws.Range("H42").Value = 14999.5
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = $null
$ws.Range("H61").Value = 3203.5881
$ws.Range("I61").Value = 2995.875
$ws.Range("K61").Value = 2995.875
$ws.Range("M61").Value = -2783.875
$ws.Range("H97").Value = 1181.8125
$ws.Range("I97").Value = 1177.2667
$ws.Range("K97").Value = 1177.2667
$ws.Range("M97").Value = -681.2666999999999
$ws.Range("H110").Value = 1770.7059
$ws.Range("J110").Value = 2400
$ws.Range("L110").Value = 2400
$ws.Range("N110").Value = -6490
$ws.Range("H132").Value = 8888.777
$ws.Range("I132").Value = 1999.5
$ws.Range("K132").Value = 5998.5
$ws.Range("M132").Value = -3468.5
$ws.Range("H136").Value = 3203.5881
$ws.Range("I136").Value = 2995.875
$ws.Range("K136").Value = 8987.625
$ws.Range("M136").Value = -6437.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I107").Value = 2482692.5
$ws.Range("J107").Value = 1463.3334
$ws.Range("K107").Value = 2482692.5
$ws.Range("L107").Value = 1463.3334
$ws.Range("M107").Value = -2480772.5
$ws.Range("N107").Value = -5303.3334
$ws.Range("H134").Value = 2514.6
$ws.Range("I134").Value = 1544.1428
$ws.Range("J134").Value = 3363.75
$ws.Range("K134").Value = 4632.428400000001
$ws.Range("L134").Value = 10091.25
$ws.Range("M134").Value = -2097.428400000001
$ws.Range("N134").Value = -15161.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3476250.8
$ws.Range("I31").Value = 2669.9524
$ws.Range("J31").Value = 8339263.5
$ws.Range("K31").Value = 2669.9524
$ws.Range("L31").Value = 8339263.5
$ws.Range("M31").Value = -2374.9524
$ws.Range("N31").Value = -8339853.5
$ws.Range("H34").Value = 3476250.8
$ws.Range("I34").Value = 2669.9524
$ws.Range("J34").Value = 8339263.5
$ws.Range("K34").Value = 2669.9524
$ws.Range("L34").Value = 8339263.5
$ws.Range("M34").Value = -2467.9524
$ws.Range("N34").Value = -8339667.5
$ws.Range("H58").Value = 2209.1785
$ws.Range("I58").Value = 1806.2106
$ws.Range("J58").Value = 3059.889
$ws.Range("K58").Value = 1806.2106
$ws.Range("L58").Value = 3059.889
$ws.Range("M58").Value = -1603.2106
$ws.Range("N58").Value = -3465.889
$ws.Range("H99").Value = 3716.9285
$ws.Range("I99").Value = 3437.4285
$ws.Range("K99").Value = 3437.4285
$ws.Range("M99").Value = -1939.4285
$ws.Range("H126").Value = 3716.9285
$ws.Range("I126").Value = 3437.4285
$ws.Range("K126").Value = 10312.2855
$ws.Range("M126").Value = -7842.2855
$ws.Range("H134").Value = 4004.9524
$ws.Range("I134").Value = 4004.9524
$ws.Range("K134").Value = 12014.8572
$ws.Range("M134").Value = -9479.8572
$ws.Range("H136").Value = 2209.1785
$ws.Range("I136").Value = 1806.2106
$ws.Range("J136").Value = 3059.889
$ws.Range("K136").Value = 5418.6318
$ws.Range("L136").Value = 9179.667000000001
$ws.Range("M136").Value = -2868.6318
$ws.Range("N136").Value = -14279.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 48.5
$ws.Range("I7").Value = 95
$ws.Range("K7").Value = 285
$ws.Range("M7").Value = -173
$ws.Range("H13").Value = 196.6
$ws.Range("I13").Value = 195.75
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 587.25
$ws.Range("L13").Value = 600
$ws.Range("M13").Value = -419.25
$ws.Range("N13").Value = -936
$ws.Range("H32").Value = 166671250
$ws.Range("J32").Value = 27783128
$ws.Range("L32").Value = 83349384
$ws.Range("N32").Value = -83349950
$ws.Range("H115").Value = 304725
$ws.Range("I115").Value = 1125
$ws.Range("J115").Value = 380625
$ws.Range("K115").Value = 3375
$ws.Range("L115").Value = 1141875
$ws.Range("M115").Value = -2200
$ws.Range("N115").Value = -1144225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6145.857
$ws.Range("I70").Value = 4753.5
$ws.Range("J70").Value = 8002.3335
$ws.Range("K70").Value = 4753.5
$ws.Range("L70").Value = 8002.3335
$ws.Range("M70").Value = -4483.5
$ws.Range("N70").Value = -8542.333500000001
$ws.Range("H73").Value = 6145.857
$ws.Range("I73").Value = 4753.5
$ws.Range("J73").Value = 8002.3335
$ws.Range("K73").Value = 4753.5
$ws.Range("L73").Value = 8002.3335
$ws.Range("M73").Value = -3817.5
$ws.Range("N73").Value = -9874.333500000001
$ws.Range("H132").Value = 34485510
$ws.Range("I132").Value = 76925336
$ws.Range("J132").Value = 3154.125
$ws.Range("K132").Value = 230776008
$ws.Range("L132").Value = 9462.375
$ws.Range("M132").Value = -230773478
$ws.Range("N132").Value = -14522.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1829.6923
$ws.Range("I7").Value = 1428.7
$ws.Range("K7").Value = 1428.7
$ws.Range("M7").Value = -1316.7
$ws.Range("H68").Value = 999.5
$ws.Range("I68").Value = 999.5
$ws.Range("K68").Value = 999.5
$ws.Range("M68").Value = -250.5
$ws.Range("H71").Value = 999.5
$ws.Range("I71").Value = 999.5
$ws.Range("K71").Value = 4997.5
$ws.Range("M71").Value = -1253.5
$ws.Range("H126").Value = 1829.6923
$ws.Range("I126").Value = 1428.7
$ws.Range("K126").Value = 4286.1
$ws.Range("M126").Value = -1816.1
$ws.Range("H132").Value = 4796.36
$ws.Range("I132").Value = 4927.933
$ws.Range("J132").Value = 4599
$ws.Range("K132").Value = 14783.799
$ws.Range("L132").Value = 13797
$ws.Range("M132").Value = -12253.799
$ws.Range("N132").Value = -18857
$ws.Range("H136").Value = 2276.7104
$ws.Range("I136").Value = 1930.6786
$ws.Range("K136").Value = 5792.0358
$ws.Range("M136").Value = -3242.0358

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 514000
$ws.Range("J18").Value = 514000
$ws.Range("L18").Value = 514000
$ws.Range("N18").Value = -514346
$ws.Range("H124").Value = 500429
$ws.Range("J124").Value = 500429
$ws.Range("L124").Value = 500429
$ws.Range("N124").Value = -510249
$ws.Range("H132").Value = 2229.1365
$ws.Range("I132").Value = 1902.1578
$ws.Range("K132").Value = 5706.4734
$ws.Range("M132").Value = -3176.4734
$ws.Range("H136").Value = 2783.6667
$ws.Range("I136").Value = 2159.3215
$ws.Range("K136").Value = 6477.9645
$ws.Range("M136").Value = -3927.9645
